# Insert a new weekly price record as row 452, shifting existing rows
# 452:493 down to 453:494 (dimension grows from A1:T493 to A1:T494).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 452, pushing everything below it down.
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A452").Value = 10
$ws.Range("B452").Value = "Vega Modelo de Temuco"
$ws.Range("C452").Value = "La Araucanía"
$ws.Range("D452").Value = 45166
$ws.Range("E452").Value = 9
$ws.Range("F452").Value = "Fruta"
$ws.Range("G452").Value = 100102
$ws.Range("H452").Value = "Cítricos"
$ws.Range("I452").Value = 100102006
$ws.Range("J452").Value = "Pomelo"
$ws.Range("K452").Value = "Start Ruby"
$ws.Range("L452").Value = "Primera"
$ws.Range("M452").Value = 35
$ws.Range("N452").Value = 15000
$ws.Range("O452").Value = 15000
$ws.Range("P452").Value = 15000
$ws.Range("Q452").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R452").Value = "Región de O'Higgins"
$ws.Range("S452").Value = 1000
$ws.Range("T452").Value = 15
